# 🔥 FINAL FIX: Admin Skills Assessment
# Rebuild the "Skills Assessment" sheet with the expanded column layout
# (question_id..prerequisite_skills) and a single sample question row,
# replacing the old 2-question / 11-column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - wipe everything currently on the sheet.
$ws.Cells.Clear()

# ---- Header row (row 1) ----
$headers = @(
    "question_id",
    "career_title",
    "skill_name",
    "skill_category",
    "skill_importance",
    "question_text",
    "question_type",
    "difficulty_level",
    "option_1",
    "option_2",
    "option_3",
    "option_4",
    "correct_answer",
    "score",
    "explanation",
    "course_link",
    "course_title",
    "learning_resource",
    "estimated_time",
    "prerequisite_skills"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data row (row 2) ----
$row2 = @(
    "Q001",
    "Full Stack Developer",
    "JavaScript",
    "Technical",
    "Critical",
    "What is closure in JavaScript?",
    "single",
    "Intermediate",
    "Function with access to parent scope",
    "Loop structure",
    "Data type",
    "Operator",
    "1",
    5,
    "Closure allows functions to access variables from outer scope",
    "/courses/javascript-advanced",
    "JavaScript Advanced Course",
    "MDN Web Docs",
    "40",
    "JavaScript Basics"
)

# Columns whose values look numeric but must stay TEXT (correct_answer = "1",
# estimated_time = "40"): use a leading apostrophe to force text entry (the
# same trick a person typing into Excel would use), then strip the resulting
# "quote prefix" cell style so the cell is left as a plain text value.
$textForceCols = @(13, 19)

for ($i = 0; $i -lt $row2.Count; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(2, $col)
    $value = $row2[$i]
    if ($value -is [string] -and ($textForceCols -contains $col)) {
        $cell.Value = "'" + $value
        $cell.ClearFormats()
    }
    else {
        $cell.Value = $value
    }
}
